$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "NamedLocation" worksheet right after "Collection" (so it
#    becomes the 9th sheet, pushing Work/WorkClosing/WorkCreation/WorkOpening/
#    License/RightsStatement down by one position each).
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("Collection")
$namedLocation = $wb.Worksheets.Add($null, $afterSheet)
$namedLocation.Name = "NamedLocation"

$namedLocation.Range("A1").Value = "@id"
$namedLocation.Range("B1").Value = "wgs:lat"
$namedLocation.Range("C1").Value = "wgs:long"

$rows = @(
    @("http://example.com/collection0/work0Location", 42.728104, -73.68757600000001),
    @("http://example.com/collection0/work1Location", 42.728104, -73.68757600000001),
    @("http://example.com/collection0/work2Location", 42.728104, -73.68757600000001),
    @("http://example.com/collection0/work3Location", 42.728104, -73.68757600000001),
    @("http://example.com/collection1/work4Location", 42.728104, -73.68757600000001),
    @("http://example.com/collection1/work5Location", 42.728104, -73.68757600000001),
    @("http://example.com/collection1/work6Location", 42.728104, -73.68757600000001),
    @("http://example.com/collection1/work7Location", 42.728104, -73.68757600000001),
    @("http://example.com/freestandingwork8Location", 42.728104, -73.68757600000001),
    @("http://example.com/freestandingwork9Location", 42.728104, -73.68757600000001),
    @("http://example.com/freestandingwork10Location", 42.728104, -73.68757600000001),
    @("http://example.com/freestandingwork11Location", 42.728104, -73.68757600000001)
)

$rowIndex = 2
foreach ($row in $rows) {
    $namedLocation.Cells.Item($rowIndex, 1).Value = $row[0]
    $namedLocation.Cells.Item($rowIndex, 2).Value = $row[1]
    $namedLocation.Cells.Item($rowIndex, 3).Value = $row[2]
    $rowIndex++
}

# ---------------------------------------------------------------------------
# 2. Person sheet: only inline anonymous models in a model's named graph --
#    rows 2 and 4's `relation` link moves from the generic Wikidata entity
#    to the more specific Wikipedia article (matching rows 5 and 6 already
#    in the sheet).
# ---------------------------------------------------------------------------
$person = $wb.Worksheets.Item("Person")
$person.Range("F2").Value = "http://en.wikipedia.org/wiki/Alan_Turing"
$person.Range("F4").Value = "http://en.wikipedia.org/wiki/Alan_Turing"

# ---------------------------------------------------------------------------
# 3. RightsStatement sheet: reword the `note` cell (E2).
# ---------------------------------------------------------------------------
$rights = $wb.Worksheets.Item("RightsStatement")
$rights.Range("E2").Value = "You may need to obtain other permissions for your intended use. For example, other rights such as publicity, privacy or moral rights may limit how you may use the material."
